# Oleksiy removed old Chris code
#
# The notes pages of several slides still carried a leftover hyperlink
# run (a bare CNN article URL) that nothing in the deck pointed to any
# more. Clear that stray run out of the notes body placeholder on each
# slide that still has it.

$p = $ppt.ActivePresentation

$staleUrl = "https://www.cnn.com/2020/07/28/business/starbucks-earnings-coronavirus/index.html"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $notesPage = $slide.NotesPage

    for ($j = 1; $j -le $notesPage.Shapes.Count; $j++) {
        $shape = $notesPage.Shapes.Item($j)

        if (-not $shape.HasTextFrame) {
            continue
        }

        $textRange = $shape.TextFrame.TextRange

        if ($textRange.Text -eq $staleUrl) {
            $textRange.Text = ""
        }
    }
}
